$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# Sheet "Foreign Currencies": update daily-rate derived values, and mark the
# two USD lots that came from dividend payments (not an actual FX purchase)
# as not taxable -- comment text + gain set to 0.
# ----------------------------------------------------------------------------
$fc = $wb.Worksheets.Item("Foreign Currencies")

# Row 2 (USD, bought 2019-09-05): Buy quantity EUR recalculated.
$fc.Range("B2").Value = 1247.91

# Row 3 (USD, bought 2022-04-01): FOREX not actually acquired (it came from a
# dividend payment), so no taxable gain.
$fc.Range("G3").Value = 0
$fc.Range("H3").Value = "FOREX not acquired (e.g. received dividend payments), thus gains not taxed."

# Row 4 (USD, bought 2022-07-01): same situation as row 3.
$fc.Range("G4").Value = 0
$fc.Range("H4").Value = "FOREX not acquired (e.g. received dividend payments), thus gains not taxed."

# Row 5 (USD, bought 2022-09-05 / sold 2022-10-12)
$fc.Range("B5").Value = 2567.09
$fc.Range("G5").Value = 57.06

# Row 6 (USD, bought 2022-09-05 / sold 2022-12-01)
$fc.Range("B6").Value = 849.87
$fc.Range("G6").Value = -43.76

# Row 7 (USD, bought 2022-09-22 / sold 2022-12-01)
$fc.Range("B7").Value = 135.13
$fc.Range("G7").Value = -7.45

# Summary rows
$fc.Range("G9").Value = 5.85
$fc.Range("G10").Value = 57.06
$fc.Range("G11").Value = -51.21

# ----------------------------------------------------------------------------
# Sheet "ELSTER - Summary": the forex gain/loss line (Anlage SO) mirrors the
# "Gains (incl. losses)" total from the Foreign Currencies sheet.
# ----------------------------------------------------------------------------
$elster = $wb.Worksheets.Item("ELSTER - Summary")
$elster.Range("C7").Value = 5.85
